# Resale numbers update: append the 2024-01-02 18:54 reading as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7

# Date/Time/Weekday/Week are stored as plain text in this sheet (e.g. "00"
# for the Week column), so force text formatting before writing them -
# otherwise Excel would auto-convert "2024-01-02" to a date serial and
# "00" to the number 0. Clear the formatting again afterwards so the new
# row ends up with the same (default) style as the existing data rows.
$textRange = $ws.Range("A" + $row + ":D" + $row)
$textRange.NumberFormat = "@"

$ws.Range("A$row").Value = "2024-01-02"
$ws.Range("B$row").Value = "18:54:12"
$ws.Range("C$row").Value = "Tuesday"
$ws.Range("D$row").Value = "00"

$ws.Range("E$row").Value = 140230
$ws.Range("F$row").Value = 142879
$ws.Range("G$row").Value = 171716
$ws.Range("H$row").Value = 145997
$ws.Range("I$row").Value = -1
$ws.Range("J$row").Value = 117134
$ws.Range("K$row").Value = 223780
$ws.Range("L$row").Value = 247983
$ws.Range("M$row").Value = 183558
$ws.Range("N$row").Value = 109830
$ws.Range("O$row").Value = 39719
$ws.Range("P$row").Value = 30742
$ws.Range("Q$row").Value = 71914
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 41612
$ws.Range("T$row").Value = -1

$textRange.ClearFormats()
